$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Neurology" to "Session"
$ws.Name = "Session"

# Append a new log entry as row 42, matching the text-formatted layout
# used by every other row in this log (Student ID, Subject, Log Date,
# Log Time, Type, User). Setting the number format to Text ("@") first
# keeps the numeric-looking Student ID stored as text rather than being
# auto-converted to a number.
$newRowRange = $ws.Range("A42:F42")
$newRowRange.NumberFormat = "@"

$ws.Range("A42").Value = "212024"
$ws.Range("B42").Value = "Neurology"
$ws.Range("C42").Value = "16/12/2025"
$ws.Range("D42").Value = "11:12:55"
$ws.Range("E42").Value = "Scan"
$ws.Range("F42").Value = "emp17.farah.a.youssef@gmail.com"
